$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.1365980374615512
$ws.Range("C2").Value = 0.5567307316750884
$ws.Range("D2").Value = 0.6103753252590118
$ws.Range("E2").Value = 0.7812652080177459
$ws.Range("F2").Value = 0.7768852086692588
$ws.Range("G2").Value = 51
$ws.Range("B3").Value = -0.1983056728170303
$ws.Range("C3").Value = 0.5132083419313908
$ws.Range("D3").Value = 0.4400880119005959
$ws.Range("E3").Value = 0.6633912962201086
$ws.Range("F3").Value = 0.6394855011274496
$ws.Range("G3").Value = 50
$ws.Range("B4").Value = -0.08011458471284517
$ws.Range("C4").Value = 0.50764697151692
$ws.Range("D4").Value = 0.518069484399113
$ws.Range("E4").Value = 0.7197704386810513
$ws.Range("F4").Value = 0.7227105481572865
$ws.Range("G4").Value = 49
$ws.Range("B5").Value = -0.1481904045569161
$ws.Range("C5").Value = 0.5126676976066116
$ws.Range("D5").Value = 0.4729556061377808
$ws.Range("E5").Value = 0.6877176790935222
$ws.Range("F5").Value = 0.6786684410696205
$ws.Range("G5").Value = 48
$ws.Range("B6").Value = -0.02446659572911941
$ws.Range("C6").Value = 0.4848588817992852
$ws.Range("D6").Value = 0.4654875719725148
$ws.Range("E6").Value = 0.6822664962992943
$ws.Range("F6").Value = 0.6891989838604905
$ws.Range("G6").Value = 47
$ws.Range("B7").Value = -0.1353368031892271
$ws.Range("C7").Value = 0.4984219442182552
$ws.Range("D7").Value = 0.4619043327933868
$ws.Range("E7").Value = 0.6796354410957295
$ws.Range("F7").Value = 0.6749645583256921
$ws.Range("G7").Value = 38
$ws.Range("B8").Value = -0.1049342653064903
$ws.Range("C8").Value = 0.4821091839203969
$ws.Range("D8").Value = 0.4458002174572263
$ws.Range("E8").Value = 0.6676827221496946
$ws.Range("F8").Value = 0.668480732802356
$ws.Range("G8").Value = 37
$ws.Range("B9").Value = -0.1747389245262045
$ws.Range("C9").Value = 0.508949487189903
$ws.Range("D9").Value = 0.4794270205088618
$ws.Range("E9").Value = 0.6924066872213626
$ws.Range("F9").Value = 0.6874003880098497
$ws.Range("G9").Value = 20
$ws.Range("B10").Value = -0.05712409056790357
$ws.Range("C10").Value = 0.5421730072189278
$ws.Range("D10").Value = 0.6097490892393119
$ws.Range("E10").Value = 0.7808643219147049
$ws.Range("F10").Value = 0.8105716633807137
$ws.Range("G10").Value = 13
$ws.Range("B11").Value = -0.5089197260409551
$ws.Range("C11").Value = 0.5296605432532845
$ws.Range("D11").Value = 0.5444046525602212
$ws.Range("E11").Value = 0.7378378226685192
$ws.Range("F11").Value = 0.597291140281082
